$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and Report date range) ---
$ws.Range("A8").Value = "Volume 29   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/21/2022  Through  11/27/2022"

# --- Style-changing cells: text<->number conversions ---
# Donor cells stable throughout the edit: A15 (text style), J15 (number style)
$ws.Range("J15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 1

$ws.Range("J15").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("F15").Value = 1

$ws.Range("J15").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = 4

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0"
$ws.Range("A15").Copy()
$ws.Range("D17").PasteSpecial(-4122)

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "***.*"
$ws.Range("A15").Copy()
$ws.Range("E17").PasteSpecial(-4122)

$ws.Range("J15").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 1

$ws.Range("J15").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("F22").Value = 1

$ws.Range("J15").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("C26").Value = 1

$ws.Range("J15").Copy()
$ws.Range("F26").PasteSpecial(-4122)
$ws.Range("F26").Value = 1

$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("A15").Copy()
$ws.Range("C27").PasteSpecial(-4122)

# --- Row 15 ---
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 17
$ws.Range("K15").Value = -15
$ws.Range("L15").Value = 21.428571428571
$ws.Range("M15").Value = 54.545454545454
$ws.Range("N15").Value = -29.166666666666

# --- Row 16 ---
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 37.5
$ws.Range("I16").Value = 108
$ws.Range("J16").Value = 98
$ws.Range("K16").Value = 10.204081632653
$ws.Range("L16").Value = 11.340206185567
$ws.Range("M16").Value = -30.769230769230
$ws.Range("N16").Value = -84.482758620689

# --- Row 17 ---
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 25
$ws.Range("I17").Value = 169
$ws.Range("K17").Value = 4.968944099378
$ws.Range("L17").Value = 33.070866141732
$ws.Range("M17").Value = 42.016806722689
$ws.Range("N17").Value = -40.282685512367

# --- Row 18 ---
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 13
$ws.Range("H18").Value = -7.142857142857
$ws.Range("I18").Value = 174
$ws.Range("J18").Value = 133
$ws.Range("K18").Value = 30.827067669172
$ws.Range("L18").Value = 31.818181818181
$ws.Range("M18").Value = -33.840304182509
$ws.Range("N18").Value = -88.298587760591

# --- Row 19 ---
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -31.25
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 49
$ws.Range("H19").Value = -10.204081632653
$ws.Range("I19").Value = 667
$ws.Range("J19").Value = 522
$ws.Range("K19").Value = 27.777777777777
$ws.Range("L19").Value = 57.683215130023
$ws.Range("M19").Value = 66.334164588528
$ws.Range("N19").Value = -6.582633053221

# --- Row 20 ---
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 10
$ws.Range("I20").Value = 136
$ws.Range("J20").Value = 101
$ws.Range("K20").Value = 34.653465346534
$ws.Range("L20").Value = 49.450549450549
$ws.Range("M20").Value = -6.849315068493
$ws.Range("N20").Value = -91.762568140520

# --- Row 21 ---
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -4.545454545454
$ws.Range("F21").Value = 100
$ws.Range("G21").Value = 99
$ws.Range("H21").Value = 1.010101010101
$ws.Range("I21").Value = 1276
$ws.Range("J21").Value = 1037
$ws.Range("K21").Value = 23.047251687560
$ws.Range("L21").Value = 44.018058690744
$ws.Range("M21").Value = 15.789473684210
$ws.Range("N21").Value = -73.766447368421

# --- Row 22 ---
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -50
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -80
$ws.Range("I22").Value = 12
$ws.Range("J22").Value = 16
$ws.Range("K22").Value = -25
$ws.Range("L22").Value = 9.090909090909
$ws.Range("M22").Value = -53.846153846153

# --- Row 24 ---
$ws.Range("C24").Value = 37
$ws.Range("D24").Value = 52
$ws.Range("E24").Value = -28.846153846153
$ws.Range("G24").Value = 130
$ws.Range("H24").Value = 33.076923076923
$ws.Range("I24").Value = 1728
$ws.Range("J24").Value = 1135
$ws.Range("K24").Value = 52.246696035242
$ws.Range("L24").Value = 57.233848953594
$ws.Range("M24").Value = 91.574279379157

# --- Row 25 ---
$ws.Range("C25").Value = 8
$ws.Range("E25").Value = -11.111111111111
$ws.Range("G25").Value = 50
$ws.Range("H25").Value = -16
$ws.Range("I25").Value = 422
$ws.Range("J25").Value = 410
$ws.Range("K25").Value = 2.926829268292
$ws.Range("L25").Value = 35.256410256410
$ws.Range("M25").Value = 4.714640198511

# --- Row 26 ---
$ws.Range("G26").Value = 6
$ws.Range("H26").Value = -83.333333333333
$ws.Range("I26").Value = 22
$ws.Range("K26").Value = -21.428571428571
$ws.Range("L26").Value = -12

# --- Row 27 ---
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = -40
$ws.Range("J27").Value = 57
$ws.Range("K27").Value = -15.789473684210
$ws.Range("L27").Value = 41.176470588235

